$d = $word.ActiveDocument

function Apply-ParagraphXml($findText, $bodyXml) {
    $rng = $d.Content
    $found = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Find failed for: $findText"
    }
    $rng.Expand(4) | Out-Null
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $bodyXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($xml)
}

# --- Section 1: "The statistical analysis plan (due April 28 at 11:59 pm)" -> "... (due May 5 at 11:59 pm)"
$body1 = '<w:r><w:rPr><w:rStyle w:val="Heading1Char"/></w:rPr><w:t>The statistical analysis plan</w:t></w:r>' + `
    '<w:r><w:rPr/><w:t xml:space="preserve"> (due </w:t></w:r>' + `
    '<w:r><w:rPr/><w:t xml:space="preserve">May </w:t></w:r>' + `
    '<w:r><w:rPr/><w:t>5</w:t></w:r>' + `
    '<w:r><w:rPr/><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:rPr/><w:t>at 11:59 pm</w:t></w:r>' + `
    '<w:r><w:rPr/><w:t>)</w:t></w:r>'
Apply-ParagraphXml "The statistical analysis plan (due April 28 at 11:59 pm)" $body1

# --- Section 2: "Peer review (due May 26 at 11:59 pm)" run split changes (text unchanged)
$body2 = '<w:r><w:rPr><w:rStyle w:val="Heading1Char"/></w:rPr><w:t>Peer review</w:t></w:r>' + `
    '<w:r><w:rPr><w:b w:val="1"/><w:bCs w:val="1"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:rPr/><w:t xml:space="preserve">(due May </w:t></w:r>' + `
    '<w:r><w:rPr/><w:t xml:space="preserve">26 </w:t></w:r>' + `
    '<w:r><w:rPr/><w:t>at 11:59 pm</w:t></w:r>' + `
    '<w:r><w:rPr/><w:t>)</w:t></w:r>'
Apply-ParagraphXml "Peer review (due May 26 at 11:59 pm)" $body2

# --- Section 3: "You will give a 10 minute slide presentation..." -> "You will give a 7 minute (without questions) slide presentation..."
$body3 = '<w:r><w:rPr/><w:t xml:space="preserve">You will </w:t></w:r>' + `
    '<w:r><w:rPr/><w:t>give a</w:t></w:r>' + `
    '<w:r><w:rPr/><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:rPr/><w:t>7</w:t></w:r>' + `
    '<w:r><w:rPr/><w:t xml:space="preserve"> minute</w:t></w:r>' + `
    '<w:r><w:rPr/><w:t xml:space="preserve"> (without questions)</w:t></w:r>' + `
    '<w:r><w:rPr/><w:t xml:space="preserve"> slide presentation on your project on either May 31 or June 2. </w:t></w:r>' + `
    '<w:r><w:rPr/><w:t xml:space="preserve">Presentation order will be </w:t></w:r>' + `
    '<w:r><w:rPr/><w:t>determined</w:t></w:r>' + `
    '<w:r><w:rPr/><w:t xml:space="preserve"> later in the quarter. </w:t></w:r>' + `
    '<w:r><w:rPr/><w:t>A final copy of your slides</w:t></w:r>' + `
    '<w:r><w:rPr/><w:t xml:space="preserve"> is due </w:t></w:r>' + `
    '<w:r><w:rPr/><w:t xml:space="preserve">at noon on the day you will be presenting. </w:t></w:r>' + `
    '<w:r><w:rPr/><w:t xml:space="preserve">I will provide you with an optional </w:t></w:r>' + `
    '<w:r><w:rPr/><w:t>powerpoint</w:t></w:r>' + `
    '<w:r><w:rPr/><w:t xml:space="preserve"> template for this </w:t></w:r>' + `
    '<w:r><w:rPr/><w:t>component</w:t></w:r>' + `
    '<w:r><w:rPr/><w:t xml:space="preserve"> of the project</w:t></w:r>' + `
    '<w:r><w:rPr/><w:t>, but you may use any organization you think will be effective and any program (</w:t></w:r>' + `
    '<w:r><w:rPr/><w:t>powerpoint</w:t></w:r>' + `
    '<w:r><w:rPr/><w:t xml:space="preserve">, R markdown, </w:t></w:r>' + `
    '<w:r><w:rPr/><w:t>beamer</w:t></w:r>' + `
    '<w:r><w:rPr/><w:t xml:space="preserve">, </w:t></w:r>' + `
    '<w:r><w:rPr/><w:t>etc</w:t></w:r>' + `
    '<w:r><w:rPr/><w:t>) to create your slides.</w:t></w:r>'
Apply-ParagraphXml "You will give a 10 minute slide presentation" $body3

# --- Section 4: "You will also think of 2 questions to ask about your peer review partner's presentation."
$body4 = '<w:r><w:rPr/><w:t xml:space="preserve">You will also think of </w:t></w:r>' + `
    '<w:r><w:rPr/><w:t>2</w:t></w:r>' + `
    '<w:r><w:rPr/><w:t xml:space="preserve"> question</w:t></w:r>' + `
    '<w:r><w:rPr/><w:t>s</w:t></w:r>' + `
    '<w:r><w:rPr/><w:t xml:space="preserve"> to ask about your peer review partner’s presentation.</w:t></w:r>'
Apply-ParagraphXml "You will also think of 2 questions to ask about your peer review partner" $body4
